$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.951.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.96%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.77%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -6.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.67%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2571"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.81%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06082"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07021"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.645.61"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.45%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5816"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -10.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.309"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.56"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.39%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.947.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006571"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.54%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.856.39"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.304"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.517"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.220"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "133.06"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.375"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.40"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.85%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.631"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.892"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.69%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07559"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.27%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.550"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9996"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04263"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -9.99%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.58%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5920"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9263"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.566"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8903"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +10.85%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01491"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.86%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.759"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3695"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.653"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1098"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.077"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05198"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.82%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.11%  "

# Row 50
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9997"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.04%  "

# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "28.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.55%  "
